# Apply "More case handling edits" change:
# Clean up GC data on "Peak_ID (3)" sheet - remove stray/duplicate rows so that
# Corrected Concentration logic can handle data sets ending on an even chain.

$wb = $excel.ActiveWorkbook

# --- Sheet "John_Code": move selection to E13 ---
$wsJohn = $wb.Worksheets.Item("John_Code")
$wsJohn.Activate()
$wsJohn.Range("E13").Select()

# --- Sheet "Peak_ID (3)": delete stray / duplicate rows ---
$wsPeak3 = $wb.Worksheets.Item("Peak_ID (3)")
$wsPeak3.Activate()

# Delete rows from bottom to top so row numbers of earlier deletions
# are not affected by later ones.
$wsPeak3.Rows.Item(212).Delete()
$wsPeak3.Rows.Item(200).Delete()
$wsPeak3.Rows.Item(189).Delete()
$wsPeak3.Rows.Item(103).Delete()
$wsPeak3.Rows.Item(102).Delete()

$wsPeak3.Application.ActiveWindow.ScrollRow = 171
$wsPeak3.Range("D205").Select()

# --- Sheet "Quantification w ES,IS": scroll position change ---
$wsQuantES = $wb.Worksheets.Item("Quantification w ES,IS")
$wsQuantES.Activate()
$wsQuantES.Application.ActiveWindow.ScrollRow = 148
$wsQuantES.Range("V1").Select()

# --- Re-activate the "Peak_ID (3)" tab (matches activeTab=2 in workbook.xml) ---
$wsPeak3.Activate()
